$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 9 (the last existing data row) into a new row 10. Copy/paste
# (rather than writing literal .Value assignments) preserves the text
# data-type of the numeric-looking cells ("1", "2", "3", ...), matching the
# original sheet's inlineStr-typed cells instead of coercing them to numbers.
$ws.Rows("9").Copy()
$ws.Rows("10").PasteSpecial(-4104)

# Row 10 differs from row 9: it has a "Cargador VE" entry and a numeric
# "Pajareras" count of 1 instead of "Sí".
$ws.Range("L10").Value = "RAEDIAN CARGADOR NEO 7KW SILVER"

# Write M10 as text "1" (matching the text data-type used for the other
# numeric-looking values) by copying an existing text "1" cell onto it.
$ws.Range("D10").Copy()
$ws.Range("M10").PasteSpecial(-4104)

# Row 9 previously carried two blank placeholder cells (E9, L9) because it
# was the last data row. Now that row 10 holds that placeholder pattern,
# drop the stale ones from row 9.
$ws.Range("E9").ClearContents()
$ws.Range("L9").ClearContents()
